$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 809
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = -128
$ws.Range("H2").Value = -127
$ws.Range("I2").Value = -127
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 872
$ws.Range("L2").Value = 414
$ws.Range("M2").Value = 458
$ws.Range("N2").Value = 458
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 209
$ws.Range("Q2").Value = -35
$ws.Range("R2").Value = -18
$ws.Range("S2").Value = 18
$ws.Range("T2").Value = 13
$ws.Range("U2").Value = -48
$ws.Range("V2").Value = 202
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = -15.65
$ws.Range("Y2").Value = -24.01
$ws.Range("Z2").Value = -13.68
$ws.Range("AA2").Value = 90.29000000000001
$ws.Range("AB2").Value = 135.43
$ws.Range("AC2").Value = -303
$ws.Range("AD2").Value = -12.68
$ws.Range("AE2").Value = 1139
$ws.Range("AF2").Value = 3.37
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 41803738

# Row 3
$ws.Range("D3").Value = 800
$ws.Range("E3").Value = -43
$ws.Range("F3").Value = -43
$ws.Range("G3").Value = -48
$ws.Range("H3").Value = -71
$ws.Range("I3").Value = -71
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 1171
$ws.Range("L3").Value = 219
$ws.Range("M3").Value = 952
$ws.Range("N3").Value = 952
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 258
$ws.Range("Q3").Value = -5
$ws.Range("R3").Value = -363
$ws.Range("S3").Value = 385
$ws.Range("T3").Value = 20
$ws.Range("U3").Value = -24
$ws.Range("V3").ClearContents()
$ws.Range("W3").Value = -5.35
$ws.Range("X3").Value = -8.84
$ws.Range("Y3").Value = -10.04
$ws.Range("Z3").Value = -6.93
$ws.Range("AA3").Value = 23.01
$ws.Range("AB3").Value = 284.15
$ws.Range("AC3").Value = -154
$ws.Range("AD3").Value = -85.42
$ws.Range("AE3").Value = 1904
$ws.Range("AF3").Value = 6.91
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 51534756

# Row 4
$ws.Range("D4").Value = 829
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = 20
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 1247
$ws.Range("L4").Value = 204
$ws.Range("M4").Value = 1044
$ws.Range("N4").Value = 1044
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 261
$ws.Range("Q4").Value = 40
$ws.Range("R4").Value = -108
$ws.Range("S4").Value = 71
$ws.Range("T4").Value = 30
$ws.Range("U4").Value = 10
$ws.Range("V4").ClearContents()
$ws.Range("W4").Value = 0.34
$ws.Range("X4").Value = 2.44
$ws.Range("Y4").Value = 2.03
$ws.Range("Z4").Value = 1.68
$ws.Range("AA4").Value = 19.52
$ws.Range("AB4").Value = 312.7
$ws.Range("AC4").Value = 39
$ws.Range("AD4").Value = 337.97
$ws.Range("AE4").Value = 2045
$ws.Range("AF4").Value = 6.43
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 52240638

# Row 5
$ws.Range("D5").Value = 842
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 35
$ws.Range("G5").Value = 26
$ws.Range("H5").Value = 58
$ws.Range("I5").Value = 58
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 1687
$ws.Range("L5").Value = 575
$ws.Range("M5").Value = 1112
$ws.Range("N5").Value = 1112
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 261
$ws.Range("Q5").Value = 421
$ws.Range("R5").Value = -101
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 17
$ws.Range("U5").Value = 404
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 4.2
$ws.Range("X5").Value = 6.9
$ws.Range("Y5").Value = 5.39
$ws.Range("Z5").Value = 3.96
$ws.Range("AA5").Value = 51.71
$ws.Range("AB5").Value = 336.51
$ws.Range("AC5").Value = 111
$ws.Range("AD5").Value = 219.27
$ws.Range("AE5").Value = 2173
$ws.Range("AF5").Value = 11.23
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 52240638

# Row 6
$ws.Range("D6").Value = 918
$ws.Range("E6").Value = 55
$ws.Range("F6").Value = 55
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 33
$ws.Range("I6").Value = 33
$ws.Range("K6").Value = 1623
$ws.Range("L6").Value = 454
$ws.Range("M6").Value = 1169
$ws.Range("N6").Value = 1169
$ws.Range("P6").Value = 261
$ws.Range("Q6").Value = -30
$ws.Range("R6").Value = -161
$ws.Range("S6").ClearContents()
$ws.Range("T6").Value = 19
$ws.Range("U6").Value = -49
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 5.96
$ws.Range("X6").Value = 3.59
$ws.Range("Y6").Value = 2.89
$ws.Range("Z6").Value = 1.99
$ws.Range("AA6").Value = 38.82
$ws.Range("AB6").Value = 348.46
$ws.Range("AC6").Value = 63
$ws.Range("AD6").Value = 573.05
$ws.Range("AE6").Value = 2285
$ws.Range("AF6").Value = 15.85
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 52240638

# Row 7
$ws.Range("D7").Value = 1076
$ws.Range("E7").Value = 160
$ws.Range("G7").Value = 161
$ws.Range("H7").Value = 147
$ws.Range("I7").Value = 147
$ws.Range("K7").Value = 1824
$ws.Range("L7").Value = 502
$ws.Range("M7").Value = 1321
$ws.Range("N7").Value = 1327
$ws.Range("P7").Value = 261
$ws.Range("Q7").Value = 146
$ws.Range("R7").Value = -109
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 27
$ws.Range("U7").Value = 136
$ws.Range("W7").Value = 14.84
$ws.Range("X7").Value = 13.65
$ws.Range("Y7").Value = 11.77
$ws.Range("Z7").Value = 8.52
$ws.Range("AA7").Value = 38.04
$ws.Range("AC7").Value = 281
$ws.Range("AD7").Value = 133.42
$ws.Range("AE7").Value = 2593
$ws.Range("AF7").Value = 14.46
$ws.Range("AG7").Value = 0
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 1205
$ws.Range("E8").Value = 211
$ws.Range("G8").Value = 228
$ws.Range("H8").Value = 204
$ws.Range("I8").Value = 204
$ws.Range("K8").Value = 2074
$ws.Range("L8").Value = 499
$ws.Range("M8").Value = 1575
$ws.Range("N8").Value = 1566
$ws.Range("P8").Value = 261
$ws.Range("Q8").Value = 224
$ws.Range("R8").Value = -49
$ws.Range("S8").Value = -1
$ws.Range("T8").Value = 41
$ws.Range("U8").Value = 173
$ws.Range("W8").Value = 17.53
$ws.Range("X8").Value = 16.96
$ws.Range("Y8").Value = 14.13
$ws.Range("Z8").Value = 10.49
$ws.Range("AA8").Value = 31.72
$ws.Range("AC8").Value = 391
$ws.Range("AD8").Value = 65.05
$ws.Range("AE8").Value = 3061
$ws.Range("AF8").Value = 8.31
$ws.Range("AG8").Value = 0
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 1289
$ws.Range("E9").Value = 221
$ws.Range("G9").Value = 238
$ws.Range("H9").Value = 208
$ws.Range("I9").Value = 208
$ws.Range("K9").Value = 2338
$ws.Range("L9").Value = 550
$ws.Range("M9").Value = 1788
$ws.Range("N9").Value = 1781
$ws.Range("P9").Value = 261
$ws.Range("Q9").Value = 262
$ws.Range("R9").Value = -62
$ws.Range("S9").Value = -1
$ws.Range("T9").Value = 37
$ws.Range("U9").Value = 204
$ws.Range("W9").Value = 17.18
$ws.Range("X9").Value = 16.12
$ws.Range("Y9").Value = 12.41
$ws.Range("Z9").Value = 9.42
$ws.Range("AA9").Value = 30.79
$ws.Range("AC9").Value = 398
$ws.Range("AD9").Value = 63.98
$ws.Range("AE9").Value = 3481
$ws.Range("AF9").Value = 7.31
$ws.Range("AG9").Value = 0
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
